# Turn the "Once Upon A Time" / "By Dorothy Day" heading block into a
# pandoc-style title block: a Title-styled paragraph with the title text
# split word-by-word into runs, followed by an Authors-styled paragraph
# with the author name (no "By " prefix, no bold) likewise split into
# runs. The two paragraphs in the source are wrapped in a malformed
# bookmark pair (non-numeric w:id="once-upon-a-time") that the Word
# object model does not expose via $d.Bookmarks, so it can't be removed
# through the Bookmarks collection - instead we delete the two
# paragraphs outright (which leaves the orphaned bookmark markers
# collapsed together at the front of the document) and then soak up
# each marker with an empty InsertXML before inserting the replacement
# paragraphs.

$d = $word.ActiveDocument

$paraOneStart = $d.Paragraphs(1).Range.Start
$paraTwoEnd = $d.Paragraphs(2).Range.End

$old = $d.Range($paraOneStart, $paraTwoEnd)
$old.Delete()

# The bookmarkStart/bookmarkEnd pair that used to wrap paragraph 1 is not
# removed by the delete above (it isn't a "real" bookmark as far as the
# object model is concerned) - it simply collapses to an empty pair at
# the very start of the document. Soak up both zero-width markers with
# a pair of empty InsertXML calls before inserting the new content.
$d.Range(0, 0).InsertXML("")
$d.Range(0, 0).InsertXML("")

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$titleXml = "<w:p $wNs>" +
    "<w:pPr><w:pStyle w:val='Title'/></w:pPr>" +
    "<w:r><w:t xml:space='preserve'>Once</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>Upon</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>A</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>Time</w:t></w:r>" +
    "</w:p>"

$authorsXml = "<w:p $wNs>" +
    "<w:pPr><w:pStyle w:val='Authors'/></w:pPr>" +
    "<w:r><w:t xml:space='preserve'>Dorothy</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>Day</w:t></w:r>" +
    "</w:p>"

$d.Range(0, 0).InsertXML($titleXml + $authorsXml)
